# Update weekly Jengibre (ginger) price records for Vega Monumental Concepción.
# Reassigns the Fecha (D), Volumen (J), Precio minimo (K), Precio maximo (L),
# Precio promedio ponderado (M) and Precio $/Kg (P) values for rows 2-31 to
# reflect the latest weekly data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


$ws.Cells.Item(2, 4).Value = 44488
$ws.Cells.Item(2, 10).Value = 40
$ws.Cells.Item(2, 11).Value = 16000
$ws.Cells.Item(2, 12).Value = 17000
$ws.Cells.Item(2, 13).Value = 16500
$ws.Cells.Item(2, 16).Value = 1269

$ws.Cells.Item(3, 4).Value = 44355
$ws.Cells.Item(3, 10).Value = 60
$ws.Cells.Item(3, 11).Value = 18000
$ws.Cells.Item(3, 12).Value = 20000
$ws.Cells.Item(3, 13).Value = 19000
$ws.Cells.Item(3, 16).Value = 1462

$ws.Cells.Item(4, 4).Value = 44264
$ws.Cells.Item(4, 10).Value = 40
$ws.Cells.Item(4, 11).Value = 30000
$ws.Cells.Item(4, 12).Value = 32000
$ws.Cells.Item(4, 13).Value = 31000
$ws.Cells.Item(4, 16).Value = 2385

$ws.Cells.Item(5, 4).Value = 44159
$ws.Cells.Item(5, 10).Value = 60
$ws.Cells.Item(5, 11).Value = 30000
$ws.Cells.Item(5, 12).Value = 32000
$ws.Cells.Item(5, 13).Value = 31000
$ws.Cells.Item(5, 16).Value = 2385

$ws.Cells.Item(6, 4).Value = 44691
$ws.Cells.Item(6, 10).Value = 100
$ws.Cells.Item(6, 11).Value = 12000
$ws.Cells.Item(6, 12).Value = 13000
$ws.Cells.Item(6, 13).Value = 12500
$ws.Cells.Item(6, 16).Value = 962

$ws.Cells.Item(7, 4).Value = 44610
$ws.Cells.Item(7, 10).Value = 50
$ws.Cells.Item(7, 11).Value = 17000
$ws.Cells.Item(7, 12).Value = 18000
$ws.Cells.Item(7, 13).Value = 17400
$ws.Cells.Item(7, 16).Value = 1338

$ws.Cells.Item(8, 4).Value = 44510
$ws.Cells.Item(8, 10).Value = 40
$ws.Cells.Item(8, 11).Value = 15000
$ws.Cells.Item(8, 12).Value = 16000
$ws.Cells.Item(8, 13).Value = 15500
$ws.Cells.Item(8, 16).Value = 1192

$ws.Cells.Item(9, 4).Value = 44383
$ws.Cells.Item(9, 10).Value = 50
$ws.Cells.Item(9, 11).Value = 15000
$ws.Cells.Item(9, 12).Value = 16000
$ws.Cells.Item(9, 13).Value = 15400
$ws.Cells.Item(9, 16).Value = 1185

$ws.Cells.Item(10, 4).Value = 44327
$ws.Cells.Item(10, 10).Value = 50
$ws.Cells.Item(10, 11).Value = 24000
$ws.Cells.Item(10, 12).Value = 25000
$ws.Cells.Item(10, 13).Value = 24400
$ws.Cells.Item(10, 16).Value = 1877

$ws.Cells.Item(11, 4).Value = 44705
$ws.Cells.Item(11, 10).Value = 50
$ws.Cells.Item(11, 11).Value = 10000
$ws.Cells.Item(11, 12).Value = 11000
$ws.Cells.Item(11, 13).Value = 10400
$ws.Cells.Item(11, 16).Value = 800

$ws.Cells.Item(12, 4).Value = 44462
$ws.Cells.Item(12, 10).Value = 60
$ws.Cells.Item(12, 11).Value = 14000
$ws.Cells.Item(12, 12).Value = 15000
$ws.Cells.Item(12, 13).Value = 14500
$ws.Cells.Item(12, 16).Value = 1115

$ws.Cells.Item(13, 4).Value = 44503
$ws.Cells.Item(13, 10).Value = 35
$ws.Cells.Item(13, 11).Value = 15000
$ws.Cells.Item(13, 12).Value = 16000
$ws.Cells.Item(13, 13).Value = 15429
$ws.Cells.Item(13, 16).Value = 1187

$ws.Cells.Item(14, 4).Value = 44316
$ws.Cells.Item(14, 10).Value = 50
$ws.Cells.Item(14, 11).Value = 27000
$ws.Cells.Item(14, 12).Value = 28000
$ws.Cells.Item(14, 13).Value = 27400
$ws.Cells.Item(14, 16).Value = 2108

$ws.Cells.Item(15, 4).Value = 44509
$ws.Cells.Item(15, 10).Value = 100
$ws.Cells.Item(15, 11).Value = 15000
$ws.Cells.Item(15, 12).Value = 16000
$ws.Cells.Item(15, 13).Value = 15500
$ws.Cells.Item(15, 16).Value = 1192

$ws.Cells.Item(16, 4).Value = 44425
$ws.Cells.Item(16, 10).Value = 60
$ws.Cells.Item(16, 11).Value = 14000
$ws.Cells.Item(16, 12).Value = 15000
$ws.Cells.Item(16, 13).Value = 14500
$ws.Cells.Item(16, 16).Value = 1115

$ws.Cells.Item(17, 4).Value = 44334
$ws.Cells.Item(17, 10).Value = 50
$ws.Cells.Item(17, 11).Value = 26000
$ws.Cells.Item(17, 12).Value = 28000
$ws.Cells.Item(17, 13).Value = 27200
$ws.Cells.Item(17, 16).Value = 2092

$ws.Cells.Item(18, 4).Value = 44523
$ws.Cells.Item(18, 10).Value = 40
$ws.Cells.Item(18, 11).Value = 15000
$ws.Cells.Item(18, 12).Value = 16000
$ws.Cells.Item(18, 13).Value = 15500
$ws.Cells.Item(18, 16).Value = 1192

$ws.Cells.Item(19, 4).Value = 44362
$ws.Cells.Item(19, 10).Value = 40
$ws.Cells.Item(19, 11).Value = 15000
$ws.Cells.Item(19, 12).Value = 16000
$ws.Cells.Item(19, 13).Value = 15500
$ws.Cells.Item(19, 16).Value = 1192

$ws.Cells.Item(20, 4).Value = 44320
$ws.Cells.Item(20, 10).Value = 50
$ws.Cells.Item(20, 11).Value = 26000
$ws.Cells.Item(20, 12).Value = 28000
$ws.Cells.Item(20, 13).Value = 26800
$ws.Cells.Item(20, 16).Value = 2062

$ws.Cells.Item(21, 4).Value = 44350
$ws.Cells.Item(21, 10).Value = 40
$ws.Cells.Item(21, 11).Value = 23000
$ws.Cells.Item(21, 12).Value = 25000
$ws.Cells.Item(21, 13).Value = 24000
$ws.Cells.Item(21, 16).Value = 1846

$ws.Cells.Item(22, 4).Value = 44467
$ws.Cells.Item(22, 10).Value = 100
$ws.Cells.Item(22, 11).Value = 13000
$ws.Cells.Item(22, 12).Value = 14000
$ws.Cells.Item(22, 13).Value = 13500
$ws.Cells.Item(22, 16).Value = 1038

$ws.Cells.Item(23, 4).Value = 44308
$ws.Cells.Item(23, 10).Value = 50
$ws.Cells.Item(23, 11).Value = 26000
$ws.Cells.Item(23, 12).Value = 27000
$ws.Cells.Item(23, 13).Value = 26400
$ws.Cells.Item(23, 16).Value = 2031

$ws.Cells.Item(24, 4).Value = 44435
$ws.Cells.Item(24, 10).Value = 100
$ws.Cells.Item(24, 11).Value = 13000
$ws.Cells.Item(24, 12).Value = 14000
$ws.Cells.Item(24, 13).Value = 13500
$ws.Cells.Item(24, 16).Value = 1038

$ws.Cells.Item(25, 4).Value = 44313
$ws.Cells.Item(25, 10).Value = 50
$ws.Cells.Item(25, 11).Value = 25000
$ws.Cells.Item(25, 12).Value = 26000
$ws.Cells.Item(25, 13).Value = 25600
$ws.Cells.Item(25, 16).Value = 1969

$ws.Cells.Item(26, 4).Value = 44664
$ws.Cells.Item(26, 10).Value = 50
$ws.Cells.Item(26, 11).Value = 11000
$ws.Cells.Item(26, 12).Value = 12000
$ws.Cells.Item(26, 13).Value = 11600
$ws.Cells.Item(26, 16).Value = 892

$ws.Cells.Item(27, 4).Value = 44433
$ws.Cells.Item(27, 10).Value = 100
$ws.Cells.Item(27, 11).Value = 13000
$ws.Cells.Item(27, 12).Value = 14000
$ws.Cells.Item(27, 13).Value = 13500
$ws.Cells.Item(27, 16).Value = 1038

$ws.Cells.Item(28, 4).Value = 44708
$ws.Cells.Item(28, 10).Value = 50
$ws.Cells.Item(28, 11).Value = 13000
$ws.Cells.Item(28, 12).Value = 14000
$ws.Cells.Item(28, 13).Value = 13600
$ws.Cells.Item(28, 16).Value = 1046

$ws.Cells.Item(29, 4).Value = 44453
$ws.Cells.Item(29, 10).Value = 50
$ws.Cells.Item(29, 11).Value = 14000
$ws.Cells.Item(29, 12).Value = 15000
$ws.Cells.Item(29, 13).Value = 14600
$ws.Cells.Item(29, 16).Value = 1123

$ws.Cells.Item(30, 4).Value = 44474
$ws.Cells.Item(30, 10).Value = 40
$ws.Cells.Item(30, 11).Value = 13000
$ws.Cells.Item(30, 12).Value = 14000
$ws.Cells.Item(30, 13).Value = 13500
$ws.Cells.Item(30, 16).Value = 1038

$ws.Cells.Item(31, 4).Value = 44377
$ws.Cells.Item(31, 10).Value = 40
$ws.Cells.Item(31, 11).Value = 14000
$ws.Cells.Item(31, 12).Value = 15000
$ws.Cells.Item(31, 13).Value = 14500
$ws.Cells.Item(31, 16).Value = 1115
